# Print values from google sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @(
    "810-488",
    "773-675",
    "546-336",
    "502-940",
    "117-135",
    "786-529",
    "255-146",
    "739-612",
    "309-650",
    "753-710",
    "687-309",
    "767-823",
    "667-654",
    "521-953",
    "572-269",
    "916-614"
)

# Clear out the old data rows (2 through 21), keeping the header row intact,
# then rewrite them with the updated set of codes (rows 2 through 17).
$ws.Range("A2:B21").Clear()

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $code = $codes[$i]
    $ws.Cells.Item($row, 1).Value = "https'//scooters.taxify.eu/qr/$code"
    $ws.Cells.Item($row, 2).Value = $code
}
